$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 844
$ws1.Range("F10").Value = 2194
$ws1.Range("F12").Value = 1661
$ws1.Range("F13").Value = 2789
$ws1.Range("F15").Value = 4165
$ws1.Range("F16").Value = 352
$ws1.Range("F17").Value = 175
$ws1.Range("F19").Value = 529
$ws1.Range("F20").Value = 247
$ws1.Range("F25").Value = 4069
$ws1.Range("F27").Value = 3529
$ws1.Range("F30").Value = 511
$ws1.Range("F33").Value = 355
$ws1.Range("F34").Value = 440
$ws1.Range("F35").Value = 341

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1013

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1013
$ws4.Range("F6").Value = 844
$ws4.Range("F12").Value = 2194
$ws4.Range("F14").Value = 1661
$ws4.Range("F16").Value = 2789
$ws4.Range("F18").Value = 4165
$ws4.Range("F19").Value = 352
$ws4.Range("F20").Value = 175
$ws4.Range("F22").Value = 529
$ws4.Range("F23").Value = 247
$ws4.Range("F29").Value = 4069
$ws4.Range("F31").Value = 3529
$ws4.Range("F34").Value = 511
$ws4.Range("F37").Value = 355
$ws4.Range("F38").Value = 440
$ws4.Range("F39").Value = 341
